$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    "B2" = 0.25;
    "C2" = 0.4166666666666667;
    "P2" = 0.25;
    "S2" = 0.08333333333333333;
    "P3" = 0.8;
    "S3" = 0.2;
    "B6" = 0.1176470588235294;
    "F6" = 0.1176470588235294;
    "J6" = 0.1764705882352941;
    "S6" = 0.5882352941176471;
    "B7" = 0.05555555555555555;
    "F7" = 0.1666666666666667;
    "J7" = 0.05555555555555555;
    "S7" = 0.7222222222222222;
    "D8" = 0.02857142857142857;
    "F8" = 0.05714285714285714;
    "O8" = 0.02857142857142857;
    "Q8" = 0.05714285714285714;
    "R8" = 0.08571428571428572;
    "S8" = 0.5428571428571428;
    "B9" = 0.0625;
    "F9" = 0.0625;
    "J9" = 0.125;
    "S9" = 0.75;
    "B10" = 0.0641025641025641;
    "D10" = 0.01282051282051282;
    "F10" = 0.07692307692307693;
    "J10" = 0.07692307692307693;
    "O10" = 0.02564102564102564;
    "Q10" = 0.217948717948718;
    "R10" = 0.03846153846153846;
    "S10" = 0.4871794871794872;
    "G11" = 0.1904761904761905;
    "J11" = 0.04761904761904762;
    "K11" = 0.1904761904761905;
    "L11" = 0.5238095238095238;
    "S11" = 0.04761904761904762;
    "G12" = 0.6666666666666666;
    "J12" = 0.1666666666666667;
    "K12" = 0.08333333333333333;
    "S12" = 0.08333333333333333;
    "G13" = 0.6666666666666666;
    "J13" = 0.2222222222222222;
    "S13" = 0.1111111111111111;
    "H15" = 0.1904761904761905;
    "I15" = 0.09523809523809523;
    "J15" = 0.3333333333333333;
    "K15" = 0.1428571428571428;
    "M15" = 0.04761904761904762;
    "O15" = 0.04761904761904762;
    "S15" = 0.1428571428571428;
    "F16" = 0.1111111111111111;
    "H16" = 0.1111111111111111;
    "I16" = 0.1111111111111111;
    "J16" = 0.1111111111111111;
    "M16" = 0.2222222222222222;
    "O16" = 0.1111111111111111;
    "S16" = 0.2222222222222222;
    "H17" = 0.1578947368421053;
    "I17" = 0.1052631578947368;
    "J17" = 0.5789473684210527;
    "M17" = 0.05263157894736842;
    "S17" = 0.1052631578947368;
    "I18" = 0.3333333333333333;
    "J18" = 0.1666666666666667;
    "M18" = 0.1666666666666667;
    "O18" = 0.1666666666666667;
    "S18" = 0.1666666666666667;
    "F19" = 0.008474576271186441;
    "H19" = 0.2457627118644068;
    "I19" = 0.07627118644067797;
    "J19" = 0.2966101694915254;
    "K19" = 0.1101694915254237;
    "M19" = 0.03389830508474576;
    "N19" = 0.008474576271186441;
    "O19" = 0.1016949152542373;
    "S19" = 0.1186440677966102
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}
